$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells that look numeric so Excel doesn't
# auto-convert strings like '1.0000' or '240.76' into numbers.

$ws.Range("D2").Value = "29.422.69"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.847.99"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.76"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6274"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07678"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.81"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07744"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "1.847.00"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.028"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001080"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6802"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.40"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.173"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "29.455.54"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.02"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.39"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.413"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.11"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1374"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.396"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.67"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.341"
$ws.Range("E28").Value = "  +5.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.461"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05662"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.119"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.025"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.841"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7007"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.581"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.773"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "1.228.83"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01788"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.534"
$ws.Range("E40").Value = "  +4.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9091"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.992.21"
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.77"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.89"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000120"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.150"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4015"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.977"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1149"
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.671"
$ws.Range("E51").Value = "  +0.39%  "
